$d = $word.ActiveDocument

# --- 1. Remove the pre-existing "_GoBack" bookmark (it currently sits
#        right after the GPA line near the end of the resume). Word
#        re-drops a "_GoBack" bookmark at the position of the most
#        recent edit, so doing this first - before we touch the email
#        hyperlink below - keeps bookmark names unique at every step. ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Find the "mailto:" hyperlink wrapping the email address and
#        remove the hyperlink, leaving plain (non-hyperlinked) text
#        with the same run formatting (sz/szCs) but no rStyle link. ---
$mailHyperlink = $null
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $candidate = $d.Hyperlinks.Item($i)
    if ($candidate.Address -like "mailto:*") {
        $mailHyperlink = $candidate
        break
    }
}

if ($mailHyperlink -ne $null) {
    $mailStart = $mailHyperlink.Range.Start
    $mailEnd = $mailHyperlink.Range.End
    $mailText = $mailHyperlink.Range.Text

    # Re-acquire a plain Range (not the live Hyperlink.Range) before
    # mutating, then inject the replacement run markup directly: same
    # text, same sz/szCs run formatting, but no rStyle="Hyperlink" and
    # no surrounding <w:hyperlink>. Word leaves a fresh "_GoBack"
    # bookmark right at the last edited spot, so recreate that too.
    $target = $d.Range($mailStart, $mailEnd)
    $xmlFragment = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">' + $mailText + '</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
    $target.InsertXML($xmlFragment)
}
